$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update start/finish timestamps and recompute diff (finish - start)
$ws.Range("B2").Value = 44469.51438932803
$ws.Range("C2").Value = 44469.51473591151
$ws.Range("D2").Value = 0.0003465834837962963

# Row 3: update start/finish timestamps and recompute diff (finish - start)
$ws.Range("B3").Value = 44469.51473595767
$ws.Range("C3").Value = 44469.51505631653
$ws.Range("D3").Value = 0.0003203588541666667
